# "Atualização final pós feedback"
# Mark the remaining in-progress / incomplete backlog items as COMPLETO and
# fill in their completion dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;  Date = "10/21/2025" },
    @{ Row = 6;  Date = "11/15/2025" },
    @{ Row = 7;  Date = "11/18/2025" },
    @{ Row = 8;  Date = "11/18/2025" },
    @{ Row = 10; Date = "11/06/2025" },
    @{ Row = 11; Date = "11/16/2025" },
    @{ Row = 12; Date = "11/16/2025" },
    @{ Row = 13; Date = "11/16/2025" },
    @{ Row = 14; Date = "11/15/2025" },
    @{ Row = 15; Date = "11/23/2025" }
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Cells.Item($row, 8).Value = "COMPLETO"
    $ws.Cells.Item($row, 9).Value = $u.Date
}

# The selection/view moved before saving the workbook.
$ws.Range("I16").Select()
$ws.Application.ActiveWindow.ScrollColumn = 5
